# Update the Metadata sheet of the CodeSystem-destino-fallecido workbook:
#  - Status changes from "draft" to "active"
#  - Date is bumped to the new publication timestamp
#  - Case Sensitive changes from "false" to "true"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B6").Value = "active"
$ws.Range("B8").Value = "2024-12-16T14:50:05-03:00"
$ws.Range("B17").Value = "true"
